$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell C1: "Remark" -> "Expected " --------------------------------
$ws.Range("C1").Value = "Expected "

# --- Highlight the header row (A1:C1) with bold font + yellow fill ----------
# Build the combined format on a scratch cell well outside the used range so
# the two separate format mutations (font, then fill) don't leave a stray,
# unused intermediate cell-style behind on cells that matter; then copy just
# the resulting format onto the header row and clean the scratch cell back up.
$scratch = $ws.Range("Z1")
$scratch.Font.Bold = $true
$scratch.Interior.Color = 65535
$scratch.Copy()
$ws.Range("A1:C1").PasteSpecial(-4122)
$scratch.Clear()

# --- Selection moves to E8 ----------------------------------------------------
$ws.Range("E8").Select() | Out-Null

# --- Page setup: portrait orientation ---------------------------------------
$ws.PageSetup.Orientation = 1
